$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativacao: keep "01/01/2022" as literal text, not a parsed date
$ws.Range("Z1").Formula = "=""01/01/2022"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Plain text replacements (existing cells keep their style)
$ws.Range("B10").Value = 'Fornecer aos alunos conceitos fundamentais para compreensão da Química Inorgânica por meio da experimentação, desenvolvendo a capacidade de realizarem práticas no laboratório que estimulem o seu pensamento científico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial.'
$ws.Range("C10").Value = 'Fornecer aos alunos conceitos fundamentais para compreensão da Química Inorgânica por meio da experimentação, desenvolvendo a capacidade de realizarem práticas no laboratório que estimulem o seu pensamento científico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial.'
$ws.Range("B14").Value = 'Compostos de Coordenação. Materiais inorgânicos de interesse industrial. Purificação e Identificação de Compostos Inorgânicos. Síntese de sais e obtenção de Compostos de Alumínio.'
$ws.Range("C14").Value = 'Compostos de Coordenação. Materiais inorgânicos de interesse industrial. Purificação e Identificação de Compostos Inorgânicos. Síntese de sais e obtenção de Compostos de Alumínio.'
$ws.Range("B16").Value = 'Compostos de Coordenação: Estrutura, ligações, reações e aplicações. Exemplos e aplicações de materiais inorgânicos de interesse industrial. Sínteses: Sal Simples, Sal Duplo e Sal Complexo. Preparação de Compostos de Alumínio.'
$ws.Range("C16").Value = 'Compostos de Coordenação: Estrutura, ligações, reações e aplicações. Exemplos e aplicações de materiais inorgânicos de interesse industrial. Sínteses: Sal Simples, Sal Duplo e Sal Complexo. Preparação de Compostos de Alumínio.'
$ws.Range("B19").Value = 'Serão oferecidas aulas expositivas e práticas.'
$ws.Range("C19").Value = 'Serão oferecidas aulas expositivas e práticas.'
$ws.Range("B20").Value = 'Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita.'
$ws.Range("C20").Value = 'Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita.'
$ws.Range("B21").Value = 'Será realizada uma prova escrita envolvendo o conteúdo do semestre todo.'
$ws.Range("C21").Value = 'Será realizada uma prova escrita envolvendo o conteúdo do semestre todo.'
$ws.Range("B22").Value = 'CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973.'
$ws.Range("C22").Value = 'CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973.'

# New cells: set value, then copy formatting from a sibling cell
# in the same column so they land on the correct existing style
$ws.Range("B11").Value = 'Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.'
$ws.Range("C11").Value = 'Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.'
$ws.Range("B9").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("B15").Value = 'Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Range("C15").Value = 'Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Range("B9").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B17").Value = 'Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Range("C17").Value = 'Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Range("B9").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C17").PasteSpecial(-4122)

